$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 11, shifting existing rows (and data) down by one.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new "output" / "configuration_fxe" record.
$ws.Range("A11").Value = "CHE"
$ws.Range("B11").Value = "trd_gasoline"
$ws.Range("C11").Value = "output"
$ws.Range("D11").Value = "configuration_fxe"
$ws.Range("F11").Value = "gasoline"
$ws.Range("G11").Value = 1

# Match the style used by neighbouring "configuration" rows (Parameter column).
$ws.Range("C11").Style = $ws.Range("C12").Style

# Update the active selection to mirror the authored edit.
$ws.Range("J11").Select()
